$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 67
$ws.Range("I5").Value = 99
$ws.Range("J5").Value = 51
$ws.Range("K5").Value = 99
$ws.Range("L5").Value = 51
$ws.Range("M5").Value = 13
$ws.Range("N5").Value = -275
$ws.Range("H32").Value = 5708.5264
$ws.Range("I32").Value = 5708.5264
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5708.5264
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -5421.5264

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 67
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 51
$ws.Range("K4").Value = 99
$ws.Range("L4").Value = 51
$ws.Range("M4").Value = 16
$ws.Range("N4").Value = -281

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 91.916664
$ws.Range("J7").Value = 81.333336
$ws.Range("L7").Value = 81.333336
$ws.Range("N7").Value = -307.333336
$ws.Range("H16").Value = 1674.3334
$ws.Range("I16").Value = 1654.7778
$ws.Range("K16").Value = 1654.7778
$ws.Range("M16").Value = -1367.7778
$ws.Range("H31").Value = 7144.067
$ws.Range("I31").Value = 3980.5
$ws.Range("J31").Value = 8725.85
$ws.Range("K31").Value = 3980.5
$ws.Range("L31").Value = 8725.85
$ws.Range("M31").Value = -3685.5
$ws.Range("N31").Value = -9315.85
$ws.Range("H34").Value = 7144.067
$ws.Range("I34").Value = 3980.5
$ws.Range("J34").Value = 8725.85
$ws.Range("K34").Value = 3980.5
$ws.Range("L34").Value = 8725.85
$ws.Range("M34").Value = -3778.5
$ws.Range("N34").Value = -9129.85
$ws.Range("H62").Value = 805
$ws.Range("I62").Value = 805
$ws.Range("K62").Value = 805
$ws.Range("M62").Value = -181
$ws.Range("H65").Value = 805
$ws.Range("I65").Value = 805
$ws.Range("K65").Value = 4025
$ws.Range("M65").Value = -905
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = ""
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = ""
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = 0
$ws.Range("H94").Value = 3120.5386
$ws.Range("I94").Value = 610
$ws.Range("K94").Value = 610
$ws.Range("M94").Value = -159
$ws.Range("H99").Value = 3808.3125
$ws.Range("I99").Value = 3494.5
$ws.Range("K99").Value = 3494.5
$ws.Range("M99").Value = -1996.5
$ws.Range("H113").Value = 1674.3334
$ws.Range("I113").Value = 1654.7778
$ws.Range("K113").Value = 1654.7778
$ws.Range("M113").Value = 515.2221999999999
$ws.Range("H122").Value = 1006.6875
$ws.Range("I122").Value = 966.5
$ws.Range("J122").Value = 1046.875
$ws.Range("K122").Value = 2899.5
$ws.Range("L122").Value = 3140.625
$ws.Range("M122").Value = -449.5
$ws.Range("N122").Value = -8040.625
$ws.Range("H126").Value = 3808.3125
$ws.Range("I126").Value = 3494.5
$ws.Range("K126").Value = 10483.5
$ws.Range("M126").Value = -8013.5
$ws.Range("H132").Value = 3825.8572
$ws.Range("I132").Value = 3065.2727
$ws.Range("K132").Value = 9195.8181
$ws.Range("M132").Value = -6665.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1234
$ws.Range("I3").Value = 1234
$ws.Range("K3").Value = 3702
$ws.Range("M3").Value = -3590
$ws.Range("H4").Value = 71572190
$ws.Range("I4").Value = 100200430
$ws.Range("J4").Value = 1599.5
$ws.Range("K4").Value = 300601290
$ws.Range("L4").Value = 4798.5
$ws.Range("M4").Value = -300601178
$ws.Range("N4").Value = -5022.5
$ws.Range("H8").Value = 294.5
$ws.Range("I8").Value = 294.5
$ws.Range("K8").Value = 883.5
$ws.Range("M8").Value = -744.5
$ws.Range("H34").Value = 2664.1428
$ws.Range("J34").Value = 3024.8333
$ws.Range("L34").Value = 9074.499899999999
$ws.Range("N34").Value = -9242.499899999999
$ws.Range("H39").Value = 5475
$ws.Range("J39").Value = 7033.3335
$ws.Range("L39").Value = 21100.0005
$ws.Range("N39").Value = -21688.0005
$ws.Range("H128").Value = 521420.56
$ws.Range("I128").Value = 521420.56
$ws.Range("K128").Value = 1564261.68
$ws.Range("M128").Value = -1559281.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = ""
$ws.Range("H113").Value = 5302.3
$ws.Range("I113").Value = 2304.9
$ws.Range("K113").Value = 2304.9
$ws.Range("M113").Value = -134.9000000000001
$ws.Range("H134").Value = 102608.336
$ws.Range("J134").Value = 102608.336
$ws.Range("L134").Value = 307825.008
$ws.Range("N134").Value = -312895.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -12550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 8748.5
$ws.Range("I33").Value = 9999
$ws.Range("J33").Value = 7498
$ws.Range("K33").Value = 9999
$ws.Range("L33").Value = 7498
$ws.Range("M33").Value = -9749
$ws.Range("N33").Value = -7998
$ws.Range("H36").Value = 8748.5
$ws.Range("I36").Value = 9999
$ws.Range("J36").Value = 7498
$ws.Range("K36").Value = 9999
$ws.Range("L36").Value = 7498
$ws.Range("M36").Value = -9749
$ws.Range("N36").Value = -7998
$ws.Range("H54").Value = 2283455.8
$ws.Range("J54").Value = 2283455.8
$ws.Range("L54").Value = 2283455.8
$ws.Range("N54").Value = -2284495.8
$ws.Range("H122").Value = 2060.111
$ws.Range("I122").Value = 1991.5714
$ws.Range("K122").Value = 5974.7142
$ws.Range("M122").Value = -3524.7142
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = ""
